# Update key-shortcut table: replace single-letter/punctuation key labels
# with the new function-key labels, per the commit diff.
$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "q"; New = "F4" },
    @{ Old = "w"; New = "F5" },
    @{ Old = "e"; New = "F6" },
    @{ Old = "r"; New = "F7" },
    @{ Old = "t"; New = "F8" },
    @{ Old = "y"; New = "F10" },
    @{ Old = "-"; New = "F11" },
    @{ Old = "="; New = "F12" }
)

foreach ($pair in $replacements) {
    $range = $d.Content
    $range.Find.Execute($pair.Old, $true, $true, $false, $false, $false, $true, 1, $false, $pair.New, 2)
}
